$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.819.97'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '2.532.12'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.80'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.64'
$ws.Range('E6').Value = '  +6.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.584'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.94'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.73'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').Value = '2.563.40'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.10'
$ws.Range('E16').Value = '  +6.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.868'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').Value = '42.876.18'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.02'
$ws.Range('E19').Value = '  +2.55%  '
$ws.Range('D20').Value = '0.0₃0982'
$ws.Range('E20').Value = '  -1.02%  '
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.52'
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.01'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.04'
$ws.Range('E25').Value = '  -3.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.97'
$ws.Range('E26').Value = '  -6.38%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.44'
$ws.Range('E28').Value = '  +1.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.33'
$ws.Range('E29').Value = '  +8.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.47'
$ws.Range('E30').Value = '  +4.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.12'
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.16'
$ws.Range('E32').Value = '  +2.55%  '
$ws.Range('E33').Value = '  -2.08%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0793'
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('E36').Value = '  -3.98%  '
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('E38').Value = '  +1.51%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.120'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.15'
$ws.Range('E40').Value = '  +3.00%  '
$ws.Range('B41').Value = 'ApeXProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  +3.46%  '
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.44'
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0303'
$ws.Range('E44').Value = '  -2.80%  '
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').Value = '2.041.96'
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.97'
$ws.Range('E47').Value = '  +1.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.98'
$ws.Range('E48').Value = '  -2.70%  '
$ws.Range('D49').Value = '2.779.55'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.93'
$ws.Range('E50').Value = '  -4.26%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.191'
$ws.Range('E51').Value = '  +0.07%  '
